{"js": "// Remove the centered \"Journal #1\" title paragraph from the top of the\n// document body (the rest of the document is left untouched).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text.trim() === \"Journal #1\") {\n    paragraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the centered \"Journal #1\" title paragraph from the top of the\n# document body (the rest of the document is left untouched).\n$d = $word.ActiveDocument\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Trim() -eq \"Journal #1\") {\n        $p.Range.Delete()\n    }\n}\n"}
